$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.071.77'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '1.676.04'
$ws.Range('E3').Value = '  +0.30%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'215.48"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +2.06%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = "'0.0621"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.35%  '
$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').Value = "'21.21"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.85%  '
$ws.Range('E11').Value = '  -0.80%  '
$ws.Range('D12').Value = '1.912.16'
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('D13').Value = '1.670.45'
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('E14').Value = '  +0.98%  '
$ws.Range('D15').Value = "'0.534"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.77%  '
$ws.Range('E16').Value = '  +0.87%  '
$ws.Range('D17').Value = '27.053.79'
$ws.Range('E17').Value = '  +0.48%  '
$ws.Range('E18').Value = '  +1.88%  '
$ws.Range('D19').Value = "'237.67"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.91%  '
$ws.Range('D20').Value = '0.0₃0744'
$ws.Range('E20').Value = '  +1.54%  '
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').Value = "'4.46"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.21%  '
$ws.Range('D24').Value = "'2.15"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.69%  '
$ws.Range('D25').Value = "'146.41"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.51%  '
$ws.Range('D26').Value = "'7.22"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.61%  '
$ws.Range('D27').Value = "'16.39"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.00%  '
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('D29').Value = "'1.00"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('E32').Value = '  +0.78%  '
$ws.Range('D33').Value = '1.548.97'
$ws.Range('E33').Value = '  +6.18%  '
$ws.Range('E34').Value = '  +2.02%  '
$ws.Range('D35').Value = "'1.71"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.04%  '
$ws.Range('D36').Value = "'0.600"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.72%  '
$ws.Range('E37').Value = '  -1.07%  '
$ws.Range('D38').Value = "'0.926"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.16%  '
$ws.Range('E39').Value = '  +2.11%  '
$ws.Range('E40').Value = '  +1.71%  '
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').Value = "'67.57"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.12%  '
$ws.Range('E43').Value = '  -2.49%  '
$ws.Range('E44').Value = '  -1.68%  '
$ws.Range('D45').Value = '1.820.99'
$ws.Range('E45').Value = '  +0.82%  '
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('D47').Value = "'90.58"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0108'
$ws.Range('E48').Value = '  +2.83%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = "'1.57"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.30%  '
$ws.Range('E50').Value = '  +2.69%  '
$ws.Range('D51').Value = "'8.04"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.79%  '
